# fixed harvester column in rnasamples -- holly added S.GISH to harvester in bioSamples
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "harvester" column (column B) for all data rows (2-19)
# from "Retrofitted_0641" to "S.GISH"
$ws.Range("B2:B19").Value = "S.GISH"

# Select the harvester column, matching the post-edit selection state
$ws.Columns("B:B").Select()
